$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the industrial building-type percentage breakdown text (B2) following
# local expert feedback: reduce share of S (simple stone/masonry) buildings.
$newText = "43% CR/LFM+CDM/H:1`n2% CR/LFM+CDM/H:1`n35% CR+PC/LFM+CDL/H:1`n5% S+SL/LFM+CDM/H:1`n11% S/LFM+CDM/RME/H:1`n3% S/LFBR+CDH/H:1`n1% S/LFBR+CDM/H:1"
$ws.Range("B2").Value = $newText

# Wrap the text so the multi-line breakdown displays properly, widen the
# column and increase the row height to fit the content.
$ws.Range("B2").WrapText = $true
$ws.Columns("B").ColumnWidth = 28.33
$ws.Rows(2).RowHeight = 120

# Keep the selection where a user would land after reviewing the update.
[void]$ws.Range("B4").Select()
